$d = $word.ActiveDocument

# Locate start of "Requisitos" heading paragraph
$rng1 = $d.Content
$rng1.Find.Execute("Requisitos", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPos = $rng1.Paragraphs(1).Range.Start

# Locate end of the last requisite paragraph (the one containing LOB1019)
$rng2 = $d.Content
$rng2.Find.Execute("LOB1019 -  Física II  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPos = $rng2.Paragraphs(1).Range.End

$delRange = $d.Range($startPos, $endPos)
$delRange.Delete()
